$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven list of (cell, new text value) updates derived from the commit diff.
# NumberFormat is forced to Text ("@") before the write so numeric-/percent-looking
# strings (e.g. "301.96", "-0.70%", "23") are stored as literal text, matching the
# original inline-string cell type instead of being auto-coerced to a number/percentage.
# Style is reset to "Normal" afterwards so no stray number-format style sticks to the cell.
$updates = @(
    @{Cell='D2'; Value='301.96'},
    @{Cell='E2'; Value='-0.70%'},
    @{Cell='G2'; Value='23'},
    @{Cell='D3'; Value='31.37'},
    @{Cell='E3'; Value='-2.09%'},
    @{Cell='G3'; Value='23'},
    @{Cell='D4'; Value='5.089'},
    @{Cell='E4'; Value='-2.86%'},
    @{Cell='G4'; Value='23'},
    @{Cell='D5'; Value='0.07390'},
    @{Cell='E5'; Value='-2.32%'},
    @{Cell='G5'; Value='23'},
    @{Cell='D6'; Value='2.239'},
    @{Cell='E6'; Value='46.80%'},
    @{Cell='G6'; Value='23'},
    @{Cell='D7'; Value='7.911'},
    @{Cell='E7'; Value='0.69%'},
    @{Cell='G7'; Value='23'},
    @{Cell='D8'; Value='3.816'},
    @{Cell='E8'; Value='-1.12%'},
    @{Cell='G8'; Value='23'},
    @{Cell='D9'; Value='0.9205'},
    @{Cell='E9'; Value='-0.85%'},
    @{Cell='G9'; Value='23'},
    @{Cell='D10'; Value='0.1705'},
    @{Cell='E10'; Value='0.62%'},
    @{Cell='G10'; Value='23'},
    @{Cell='D11'; Value='0.07504'},
    @{Cell='E11'; Value='-5.17%'},
    @{Cell='G11'; Value='23'},
    @{Cell='D12'; Value='0.08158'},
    @{Cell='E12'; Value='1.43%'},
    @{Cell='G12'; Value='23'},
    @{Cell='D13'; Value='0.03024'},
    @{Cell='E13'; Value='-0.02%'},
    @{Cell='G13'; Value='23'},
    @{Cell='D14'; Value='0.09930'},
    @{Cell='E14'; Value='0.29%'},
    @{Cell='G14'; Value='23'},
    @{Cell='D15'; Value='0.001492'},
    @{Cell='E15'; Value='-0.41%'},
    @{Cell='G15'; Value='23'},
    @{Cell='D16'; Value='0.006107'},
    @{Cell='E16'; Value='-2.43%'},
    @{Cell='G16'; Value='23'},
    @{Cell='D17'; Value='3.453'},
    @{Cell='E17'; Value='-0.06%'},
    @{Cell='G17'; Value='23'},
    @{Cell='D18'; Value='2.226'},
    @{Cell='E18'; Value='-0.52%'},
    @{Cell='G18'; Value='23'},
    @{Cell='D19'; Value='0.3282'},
    @{Cell='E19'; Value='-0.55%'},
    @{Cell='G19'; Value='23'},
    @{Cell='D20'; Value='0.1319'},
    @{Cell='E20'; Value='-1.01%'},
    @{Cell='G20'; Value='23'},
    @{Cell='E21'; Value='4.92%'},
    @{Cell='G21'; Value='23'},
    @{Cell='E22'; Value='0.92%'},
    @{Cell='G22'; Value='23'},
    @{Cell='D23'; Value='0.1566'},
    @{Cell='E23'; Value='-3.14%'},
    @{Cell='G23'; Value='23'},
    @{Cell='D24'; Value='0.001224'},
    @{Cell='E24'; Value='0.74%'},
    @{Cell='G24'; Value='23'},
    @{Cell='D25'; Value='0.004492'},
    @{Cell='E25'; Value='0.12%'},
    @{Cell='G25'; Value='23'},
    @{Cell='D26'; Value='0.0001297'},
    @{Cell='E26'; Value='-7.09%'},
    @{Cell='G26'; Value='23'},
    @{Cell='D27'; Value='0.0003421'},
    @{Cell='E27'; Value='92.24%'},
    @{Cell='G27'; Value='23'},
    @{Cell='G28'; Value='23'},
    @{Cell='G29'; Value='23'},
    @{Cell='G30'; Value='23'},
    @{Cell='G31'; Value='23'},
    @{Cell='G32'; Value='23'},
    @{Cell='G33'; Value='23'},
    @{Cell='G34'; Value='23'},
    @{Cell='G35'; Value='23'},
    @{Cell='G36'; Value='23'},
    @{Cell='G37'; Value='23'},
    @{Cell='G38'; Value='23'},
    @{Cell='D39'; Value='0.01742'},
    @{Cell='E39'; Value='2.09%'},
    @{Cell='G39'; Value='23'},
    @{Cell='D40'; Value='0.04511'},
    @{Cell='E40'; Value='0.50%'},
    @{Cell='G40'; Value='23'},
    @{Cell='D41'; Value='0.007354'},
    @{Cell='E41'; Value='5.46%'},
    @{Cell='G41'; Value='23'},
    @{Cell='D42'; Value='0.1353'},
    @{Cell='E42'; Value='-0.36%'},
    @{Cell='G42'; Value='23'},
    @{Cell='D43'; Value='0.002225'},
    @{Cell='E43'; Value='7.28%'},
    @{Cell='G43'; Value='23'},
    @{Cell='E44'; Value='-22.56%'},
    @{Cell='G44'; Value='23'},
    @{Cell='D45'; Value='0.00006292'},
    @{Cell='E45'; Value='2.14%'},
    @{Cell='G45'; Value='23'},
    @{Cell='E46'; Value='12.44%'},
    @{Cell='G46'; Value='23'},
    @{Cell='E47'; Value='-22.89%'},
    @{Cell='G47'; Value='23'},
    @{Cell='G48'; Value='23'},
    @{Cell='G49'; Value='23'},
    @{Cell='G50'; Value='23'},
    @{Cell='G51'; Value='23'}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
